# Update "想去人数" (F column) counts that changed between scrapes.
# Sheet "展览" (exhibition) and sheet "全部类型" (all types) both carry the
# same event list (the latter has one extra row early on), so the F-column
# row numbers differ slightly between the two sheets for the same events.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1291
$ws1.Range("F10").Value = 3544
$ws1.Range("F14").Value = 50
$ws1.Range("F16").Value = 614
$ws1.Range("F17").Value = 105
$ws1.Range("F18").Value = 771
$ws1.Range("F24").Value = 2749
$ws1.Range("F25").Value = 5239
$ws1.Range("F29").Value = 3093
$ws1.Range("F30").Value = 294
$ws1.Range("F31").Value = 2271
$ws1.Range("F33").Value = 493
$ws1.Range("F35").Value = 135
$ws1.Range("F36").Value = 184
$ws1.Range("F42").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1291
$ws4.Range("F10").Value = 3544
$ws4.Range("F15").Value = 50
$ws4.Range("F17").Value = 614
$ws4.Range("F18").Value = 105
$ws4.Range("F19").Value = 771
$ws4.Range("F25").Value = 2749
$ws4.Range("F26").Value = 5239
$ws4.Range("F30").Value = 3093
$ws4.Range("F31").Value = 294
$ws4.Range("F32").Value = 2271
$ws4.Range("F34").Value = 493
$ws4.Range("F36").Value = 135
$ws4.Range("F37").Value = 184
$ws4.Range("F43").Value = 6
